$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E").Delete()
